$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Recalculated strikeout (K) values for rows 2-21 (column G)
$kValues = @{
    2  = 1
    3  = 4
    4  = 3
    5  = 4
    6  = 1
    7  = 3
    8  = 6
    9  = 5
    10 = 2
    11 = 3
    12 = 5
    13 = 5
    14 = 4
    15 = 3
    16 = 4
    17 = 6
    18 = 4
    19 = 2
    20 = 2
    21 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
